$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D18: replace AVEDEV with MEDIAN
$ws.Range("D18").Formula = "=MEDIAN(B1:B18)"

# E18: drop formula, becomes a plain numeric value (0)
$ws.Range("E18").ClearContents()
$ws.Range("E18").Value = 0

# F18: now holds the STDEV.P formula (previously in E18)
$ws.Range("F18").Formula = "=STDEV.P(B1:B18)"

# G18: now holds the VAR.P formula (previously in F18); old ratio formula removed
$ws.Range("G18").Formula = "=VAR.P(B1:B18)"

# H18 unchanged: VAR.S(B1:B18)

# Update the active selection to D18, matching the saved view state
$ws.Range("D18").Select()
